$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Simple value edits on existing rows (dates / counts)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "22-06-2020"
$ws.Range("B5").Value = 44

# ---------------------------------------------------------------------------
# 2. Row 9 (Id) - reword the Definition / Attribute description text
#    ("...ecosystems classifications..." -> "Identification code of the
#    line.", "Numeric value for each polygon." -> "...for each line.")
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "Identification code of the line."
$ws.Range("E9").Value = "Numeric value for each line. "

# ---------------------------------------------------------------------------
# 3. Row 11 (Trail) - capitalise definition, add Type "String", and replace
#    the plain-text value list with a formatted rich-text description.
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Category of trail. "
$ws.Range("D11").Value = "String"

$rng = $ws.Range("E11")
$rng.Value = "Values:
Rogue = unknown or pedestrian-made pathway. 
Main = decided by University of Victoria officials. 
CJVI edge = walked along the edge of the CJVI property.
GO trees = walked along the Garry Oak trees. "

# run boundaries (1-based character positions) ------------------------------
# "Values:\n"                                             -> underline
# "Rogue"                                                  -> bold
# " = unknown or pedestrian-made pathway. \n"              -> plain
# "Main"                                                   -> bold
# " = decided by University of Victoria officials. \n"     -> plain
# "CJVI edge "                                             -> bold
# "= walked along the edge of the CJVI property.\n"        -> plain
# "GO trees "                                              -> bold
# "= walked along the Garry Oak trees. "                   -> plain
$rng.Characters(1, 8).Font.Underline = $true
$rng.Characters(1, 8).Font.Name = "Times New Roman"
$rng.Characters(1, 8).Font.Size = 9

$rng.Characters(9, 5).Font.Bold = $true
$rng.Characters(9, 5).Font.Name = "Times New Roman"
$rng.Characters(9, 5).Font.Size = 9

$rng.Characters(14, 40).Font.Name = "Times New Roman"
$rng.Characters(14, 40).Font.Size = 9

$rng.Characters(54, 4).Font.Bold = $true
$rng.Characters(54, 4).Font.Name = "Times New Roman"
$rng.Characters(54, 4).Font.Size = 9

$rng.Characters(58, 49).Font.Name = "Times New Roman"
$rng.Characters(58, 49).Font.Size = 9

$rng.Characters(107, 10).Font.Bold = $true
$rng.Characters(107, 10).Font.Name = "Times New Roman"
$rng.Characters(107, 10).Font.Size = 9

$rng.Characters(117, 46).Font.Name = "Times New Roman"
$rng.Characters(117, 46).Font.Size = 9

$rng.Characters(163, 9).Font.Bold = $true
$rng.Characters(163, 9).Font.Name = "Times New Roman"
$rng.Characters(163, 9).Font.Size = 9

$rng.Characters(172, 36).Font.Name = "Times New Roman"
$rng.Characters(172, 36).Font.Size = 9

# ---------------------------------------------------------------------------
# 4. New row 12 - "Comments" attribute row, appended after "Trail"
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Comments"
$ws.Range("A12").Style = "Attribute name"

$ws.Range("B12").Value = "Comments on data."
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = "String"
$ws.Range("E12").Value = ""
$ws.Range("B12:E12").Style = "metaText"

# restore word-wrap lost by the named-style assignment above (matches the
# wrapText alignment used by every other data row)
$ws.Range("A12:E12").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Sheet view tidy-up: drop the frozen "topLeftCell", move the selection
# ---------------------------------------------------------------------------
$ws.Range("E9").Select()

# ---------------------------------------------------------------------------
# 6. Page setup: portrait orientation
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
